# Metabolomics_computational_analysis.xlsx
# "update building block types in metabolights, RNA, DNA and metabolite templates"
#
# Content-level edits (the rest of the canonical-XML diff — dropped revision
# pointers, bookViews, xr:uid noise, re-ordered font children, table
# totals-row bookkeeping, etc. — is re-save churn produced by whatever tool
# wrote the "after" workbook, not an addressable cell/property edit, so it
# is left alone here and is regenerated consistently by the runtime itself).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("isa_template")
$ws2 = $wb.Worksheets.Item("4COM03_Metabolomics")

# TEMPLATE sheet: bump template version 1.1.8 -> 1.1.9
$ws1.Range("B4").Value = "1.1.9"

# Annotation table header row: rename two building-block columns.
# Excel keeps the table's tableColumn/@name in sync with the header cell
# automatically, so editing the header cells is sufficient.
$ws2.Range("K1").Value = "Component [metabolite assignment file]"
$ws2.Range("Q1").Value = "Output [Data]"

# Term Accession Number (DPBO:0000077): DPBO purl moved to the new
# nfdi4plants ontology host.
$ws2.Range("D2").Value = "http://purl.org/nfdi4plants/ontology/dpbo/DPBO_1000178"

# Parameter [metabolite assignment file] value cell ("Term Source REF
# (DPBO:0000077)" row) is cleared back to blank now that the building
# block above is a Component instead of a Parameter.
$ws2.Range("L2").Value = ""
